$wb = $excel.ActiveWorkbook

$active   = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")
$config   = $wb.Worksheets.Item("Config")

# ---------------------------------------------------------------------------
# 1. Inactive sheet: insert a new row at the top of the data (row 2) for the
#    task that just got completed ("zoom in needs to keep pixels clear...",
#    previously row 4 / Id 40 on the Active sheet). All existing Inactive
#    rows shift down by one.
# ---------------------------------------------------------------------------

# Grab style templates from the existing first data row before inserting, so
# the new row matches the plain "data row" styling (no explicit style / s="0").
$styleA = $inactive.Range("A2").Style
$styleB = $inactive.Range("B2").Style
$styleC = $inactive.Range("C2").Style
$styleD = $inactive.Range("D2").Style
$styleE = $inactive.Range("E2").Style
$styleF = $inactive.Range("F2").Style

$inactive.Rows.Item(2).Insert()

$inactive.Range("A2").Value = 40
$inactive.Range("B2").Value = "zoom in needs to keep pixels clear instead of letting it blur together"
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "Task"

# Date-like columns: format as text first so Excel stores the literal string
# instead of silently converting it to a date serial number, then restore the
# plain data-row style.
$inactive.Range("E2").NumberFormat = "@"
$inactive.Range("E2").Value = "8/13/2018"
$inactive.Range("E2").Style = $styleE

$inactive.Range("F2").NumberFormat = "@"
$inactive.Range("F2").Value = "8/15/2018"
$inactive.Range("F2").Style = $styleF

$inactive.Range("A2").Style = $styleA
$inactive.Range("B2").Style = $styleB
$inactive.Range("C2").Style = $styleC
$inactive.Range("D2").Style = $styleD

# ---------------------------------------------------------------------------
# 2. Active sheet: row 4 held that now-completed task. Replace it in place
#    with the newly discovered bug report (next Id in sequence).
# ---------------------------------------------------------------------------

$styleE4 = $active.Range("E4").Style

$active.Range("A4").Value = 41
$active.Range("B4").Value = "bug: on some colors (oranges esp.) adjusting the saturation gives a too bright color - keep it in the gray range"
$active.Range("C4").Value = "Todo"
$active.Range("D4").Value = "Task"

$active.Range("E4").NumberFormat = "@"
$active.Range("E4").Value = "8/15/2018"
$active.Range("E4").Style = $styleE4

# ---------------------------------------------------------------------------
# 3. Config sheet: bump "Max Id" now that Id 41 has been used.
# ---------------------------------------------------------------------------

$config.Range("F2").Value = 41

Write-Host "edit applied"
